# Update countries & provincias Spain
# Refreshes the COVID-19 country figures (Casos totales/Nuevos casos/
# Casos activos/Recuperados/Muertes hoy/Muertes) and the "last updated"
# banner, matching the upstream data-source refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Datos actualizados ..." banner (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 24 de Octubre de 2020 a las 22:44"

# Row 4 (Estados Unidos)
$ws.Cells.Item(4,2).Value = 8814779   # B4: 8800758 -> 8814779
$ws.Cells.Item(4,3).Value = 67751   # C4: 53805 -> 67751
$ws.Cells.Item(4,4).Value = 5733577   # D4: 5719644 -> 5733577
$ws.Cells.Item(4,5).Value = 2851217   # E4: 2851295 -> 2851217
$ws.Cells.Item(4,7).Value = 701   # G4: 535 -> 701
$ws.Cells.Item(4,8).Value = 229985   # H4: 229819 -> 229985

# Row 5 (India)
$ws.Cells.Item(5,2).Value = 7863533   # B5: 7863450 -> 7863533
$ws.Cells.Item(5,3).Value = 49865   # C5: 49782 -> 49865
$ws.Cells.Item(5,4).Value = 7075273   # D5: 7075132 -> 7075273
$ws.Cells.Item(5,5).Value = 669701   # E5: 669760 -> 669701
$ws.Cells.Item(5,7).Value = 567   # G5: 566 -> 567
$ws.Cells.Item(5,8).Value = 118559   # H5: 118558 -> 118559

# Row 20 (Alemania)
$ws.Cells.Item(20,2).Value = 427398   # B20: 424527 -> 427398
$ws.Cells.Item(20,3).Value = 10048   # C20: 7177 -> 10048
$ws.Cells.Item(20,5).Value = 103189   # E20: 100327 -> 103189
$ws.Cells.Item(20,7).Value = 19   # G20: 10 -> 19
$ws.Cells.Item(20,8).Value = 10109   # H20: 10100 -> 10109

# Row 28 (Israel)
$ws.Cells.Item(28,2).Value = 309413   # B28: 309374 -> 309413
$ws.Cells.Item(28,3).Value = 573   # C28: 534 -> 573
$ws.Cells.Item(28,4).Value = 291206   # D28: 291130 -> 291206
$ws.Cells.Item(28,5).Value = 15835   # E28: 15878 -> 15835
$ws.Cells.Item(28,7).Value = 43   # G28: 37 -> 43
$ws.Cells.Item(28,8).Value = 2372   # H28: 2366 -> 2372

# Row 33 (Canada)
$ws.Cells.Item(33,2).Value = 213959   # B33: 213881 -> 213959
$ws.Cells.Item(33,3).Value = 2227   # C33: 2149 -> 2227
$ws.Cells.Item(33,4).Value = 179636   # D33: 179621 -> 179636
$ws.Cells.Item(33,5).Value = 24401   # E33: 24338 -> 24401

# Row 51 (Costa Rica)
$ws.Cells.Item(51,2).Value = 103088   # B51: 101826 -> 103088
$ws.Cells.Item(51,3).Value = 1262   # C51: 0 -> 1262
$ws.Cells.Item(51,4).Value = 62037   # D51: 61662 -> 62037
$ws.Cells.Item(51,5).Value = 39769   # E51: 38899 -> 39769
$ws.Cells.Item(51,7).Value = 17   # G51: 0 -> 17
$ws.Cells.Item(51,8).Value = 1282   # H51: 1265 -> 1282

# Row 53 (Etiopia)
$ws.Cells.Item(53,2).Value = 92858   # B53: 92229 -> 92858
$ws.Cells.Item(53,3).Value = 629   # C53: 0 -> 629
$ws.Cells.Item(53,4).Value = 46842   # D53: 46118 -> 46842
$ws.Cells.Item(53,5).Value = 44597   # E53: 44711 -> 44597
$ws.Cells.Item(53,7).Value = 19   # G53: 0 -> 19
$ws.Cells.Item(53,8).Value = 1419   # H53: 1400 -> 1419

# Row 95 (Costa de Marfil)
$ws.Cells.Item(95,2).Value = 20429   # B95: 20405 -> 20429
$ws.Cells.Item(95,3).Value = 24   # C95: 0 -> 24
$ws.Cells.Item(95,4).Value = 20137   # D95: 20100 -> 20137
$ws.Cells.Item(95,5).Value = 171   # E95: 184 -> 171

# Row 100 (Zambia)
$ws.Cells.Item(100,2).Value = 16117   # B100: 16095 -> 16117
$ws.Cells.Item(100,3).Value = 22   # C100: 0 -> 22
$ws.Cells.Item(100,5).Value = 590   # E100: 570 -> 590
$ws.Cells.Item(100,7).Value = 2   # G100: 0 -> 2
$ws.Cells.Item(100,8).Value = 348   # H100: 346 -> 348

# Row 106 (Mozambique)
$ws.Cells.Item(106,2).Value = 11895   # B106: 11748 -> 11895
$ws.Cells.Item(106,3).Value = 147   # C106: 0 -> 147
$ws.Cells.Item(106,4).Value = 9244   # D106: 9234 -> 9244
$ws.Cells.Item(106,5).Value = 2566   # E106: 2432 -> 2566
$ws.Cells.Item(106,7).Value = 3   # G106: 0 -> 3
$ws.Cells.Item(106,8).Value = 85   # H106: 82 -> 85

# Row 108 (Maldivas)
$ws.Cells.Item(108,2).Value = 11421   # B108: 11391 -> 11421
$ws.Cells.Item(108,3).Value = 30   # C108: 0 -> 30
$ws.Cells.Item(108,4).Value = 10472   # D108: 10428 -> 10472
$ws.Cells.Item(108,5).Value = 912   # E108: 926 -> 912

# Row 129 (Republica de Yibuti)
$ws.Cells.Item(129,2).Value = 5530   # B129: 5528 -> 5530
$ws.Cells.Item(129,3).Value = 2   # C129: 0 -> 2
$ws.Cells.Item(129,4).Value = 5398   # D129: 5393 -> 5398
$ws.Cells.Item(129,5).Value = 71   # E129: 74 -> 71

# Row 130 (Trinidad yTobago)
$ws.Cells.Item(130,2).Value = 5503   # B130: 5487 -> 5503
$ws.Cells.Item(130,3).Value = 16   # C130: 0 -> 16
$ws.Cells.Item(130,4).Value = 4018   # D130: 3945 -> 4018
$ws.Cells.Item(130,5).Value = 1380   # E130: 1438 -> 1380
$ws.Cells.Item(130,7).Value = 1   # G130: 0 -> 1
$ws.Cells.Item(130,8).Value = 105   # H130: 104 -> 105

# Row 152 (Republica de Chipre)
$ws.Cells.Item(152,2).Value = 3444   # B152: 3314 -> 3444
$ws.Cells.Item(152,3).Value = 130   # C152: 0 -> 130
$ws.Cells.Item(152,5).Value = 1537   # E152: 1407 -> 1537

# Row 159 (Sierra Leona)
$ws.Cells.Item(159,2).Value = 2345   # B159: 2343 -> 2345
$ws.Cells.Item(159,3).Value = 2   # C159: 0 -> 2
$ws.Cells.Item(159,4).Value = 1784   # D159: 1782 -> 1784

# Row 162 (Yemen)
$ws.Cells.Item(162,4).Value = 1360   # D162: 1354 -> 1360
$ws.Cells.Item(162,5).Value = 101   # E162: 107 -> 101

# Rows 216/217: "Montserrat" and "Islas Malvinas" swap places in the
# country list (their shared-string/list order is swapped upstream),
# and each takes its own refreshed Casos activos (D) / Muertes (H) with it.
$ws.Cells.Item(216,1).Value = "Islas Malvinas"
$ws.Cells.Item(216,4).Value = 13   # D216: 12 -> 13
$ws.Cells.Item(216,8).Value = 0    # H216: 1 -> 0

$ws.Cells.Item(217,1).Value = "Montserrat"
$ws.Cells.Item(217,4).Value = 12   # D217: 13 -> 12
$ws.Cells.Item(217,8).Value = 1    # H217: 0 -> 1
